$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 header additions: new "Sign" column header next to each Degrees column ---
foreach ($addr in @("D3", "H3", "M3", "Q3")) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = "Sign"
}

# --- Insert a new data row (LEBL / 06L/24R) above the current row 6 (SCQ / 06R/24L) ---
$ws.Rows.Item(6).Insert()

$ws.Range("A6").Value = "LEBL"
$ws.Range("B6").Value = "06L/24R"
$ws.Range("C6").Value = "8"
$ws.Range("D6").Value = "N"
$ws.Range("E6").Value = 41
$ws.Range("F6").Value = 17.69
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = "E"
$ws.Range("I6").Value = 2
$ws.Range("J6").Value = 4.32
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 10
$ws.Range("M6").Value = "N"
$ws.Range("N6").Value = 41
$ws.Range("O6").Value = 18.34
$ws.Range("P6").Value = 1
$ws.Range("Q6").Value = "E"
$ws.Range("R6").Value = 2
$ws.Range("S6").Value = 6.22
$ws.Range("T6").Value = 0
$ws.Range("U6").Value = 45

# --- Update selection to match the committed state ---
$ws.Range("U6").Select()
